# Rename the worksheet "Embedding" -> "ATLAS" (only content change in this revision).
$wb = $excel.ActiveWorkbook

$ws = $null
try {
    $ws = $wb.Worksheets.Item("Embedding")
} catch {
    $ws = $null
}

if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$ws.Name = "ATLAS"
